# Update cryptos list: price (D) and volume(1h) (E) changes, plus
# three-row reorderings (rows 29-31 and rows 39-40) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as text so values like "596.50" or "1.00" are not
# silently renumbered (trailing zeros dropped) by Excel's auto type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '65.549.63'
$ws.Cells.Item(2, 5).Value = '  -0.36%  '
$ws.Cells.Item(3, 4).Value = '2.657.80'
$ws.Cells.Item(3, 5).Value = '  -0.74%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).Value = '596.50'
$ws.Cells.Item(5, 5).Value = '  -0.74%  '
$ws.Cells.Item(6, 4).Value = '155.73'
$ws.Cells.Item(6, 5).Value = '  -0.44%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).Value = '0.625'
$ws.Cells.Item(8, 5).Value = '  +6.55%  '
$ws.Cells.Item(9, 5).Value = '  +3.17%  '
$ws.Cells.Item(10, 4).Value = '0.395'
$ws.Cells.Item(10, 5).Value = '  -0.51%  '
$ws.Cells.Item(11, 5).Value = '  -1.70%  '
$ws.Cells.Item(12, 5).Value = '  +0.40%  '
$ws.Cells.Item(13, 4).Value = '28.78'
$ws.Cells.Item(13, 5).Value = '  -2.30%  '
$ws.Cells.Item(15, 4).Value = '3.135.30'
$ws.Cells.Item(15, 5).Value = '  -0.69%  '
$ws.Cells.Item(16, 4).Value = '65.395.66'
$ws.Cells.Item(16, 5).Value = '  -0.31%  '
$ws.Cells.Item(17, 4).Value = '2.659.36'
$ws.Cells.Item(17, 5).Value = '  -1.02%  '
$ws.Cells.Item(18, 4).Value = '12.64'
$ws.Cells.Item(18, 5).Value = '  +0.35%  '
$ws.Cells.Item(19, 4).Value = '4.76'
$ws.Cells.Item(19, 5).Value = '  -1.11%  '
$ws.Cells.Item(20, 4).Value = '7.47'
$ws.Cells.Item(20, 5).Value = '  -1.67%  '
$ws.Cells.Item(21, 4).Value = '348.80'
$ws.Cells.Item(21, 5).Value = '  -0.60%  '
$ws.Cells.Item(22, 5).Value = '  -0.01%  '
$ws.Cells.Item(23, 4).Value = '68.92'
$ws.Cells.Item(23, 5).Value = '  -1.58%  '
$ws.Cells.Item(24, 5).Value = '  +1.94%  '
$ws.Cells.Item(25, 4).Value = '9.62'
$ws.Cells.Item(25, 5).Value = '  -2.06%  '
$ws.Cells.Item(26, 4).Value = '1.66'
$ws.Cells.Item(26, 5).Value = '  +2.17%  '
$ws.Cells.Item(27, 4).Value = '1.59'
$ws.Cells.Item(27, 5).Value = '  -2.39%  '
$ws.Cells.Item(28, 5).Value = '  -2.70%  '
$ws.Cells.Item(29, 2).Value = 'Bittensor'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(29, 4).Value = '540.94'
$ws.Cells.Item(29, 5).Value = '  +2.03%  '
$ws.Cells.Item(30, 2).Value = 'Aptos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(30, 4).Value = '7.95'
$ws.Cells.Item(30, 5).Value = '  -2.57%  '
$ws.Cells.Item(31, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(31, 4).Value = '1.00'
$ws.Cells.Item(31, 5).Value = '  -0.04%  '
$ws.Cells.Item(32, 5).Value = '  -2.82%  '
$ws.Cells.Item(33, 5).Value = '  -0.20%  '
$ws.Cells.Item(34, 4).Value = '6.41'
$ws.Cells.Item(34, 5).Value = '  -2.33%  '
$ws.Cells.Item(35, 5).Value = '  +0.52%  '
$ws.Cells.Item(36, 5).Value = '  -1.55%  '
$ws.Cells.Item(37, 5).Value = '  -0.37%  '
$ws.Cells.Item(38, 5).Value = '  +0.07%  '
$ws.Cells.Item(39, 2).Value = 'Monero'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(39, 4).Value = '155.47'
$ws.Cells.Item(39, 5).Value = '  -3.31%  '
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).Value = '1.92'
$ws.Cells.Item(40, 5).Value = '  -2.38%  '
$ws.Cells.Item(41, 5).Value = '  +0.00%  '
$ws.Cells.Item(42, 4).Value = '161.18'
$ws.Cells.Item(42, 5).Value = '  -2.81%  '
$ws.Cells.Item(43, 5).Value = '  -0.57%  '
$ws.Cells.Item(44, 5).Value = '  +3.29%  '
$ws.Cells.Item(45, 5).Value = '  -2.52%  '
$ws.Cells.Item(46, 4).Value = '22.50'
$ws.Cells.Item(46, 5).Value = '  -2.46%  '
$ws.Cells.Item(47, 4).Value = '0.636'
$ws.Cells.Item(47, 5).Value = '  -2.11%  '
$ws.Cells.Item(48, 5).Value = '  -2.87%  '
$ws.Cells.Item(49, 4).Value = '0.0992'
$ws.Cells.Item(49, 5).Value = '  +0.79%  '
$ws.Cells.Item(50, 4).Value = '0.0₆0251'
$ws.Cells.Item(50, 5).Value = '  +7.18%  '
$ws.Cells.Item(51, 4).Value = '19.67'
$ws.Cells.Item(51, 5).Value = '  -3.42%  '
